$wb = $excel.ActiveWorkbook
$tmpl = $wb.Worksheets.Item("Template")
$tmpl.Copy($null, $tmpl)
$new = $wb.Worksheets.Item("Template (2)")
$new.Name = "HDD__20230427_155542"

# --- Header row ---
$new.Range("B1").Value = "HDD"
$new.Range("D1").Value = ""

# --- OPTIONS block ---
$new.Range("B4").Value = "/home/lukasbosshart/workspace/defects4j/bugs/csv_4_b/"
$new.Range("B5").Value = "src/main/java"
$new.Range("B6").Value = "src/test/java"
$new.Range("B7").Value = "org.apache.commons.csv.CSVParserTest#testNoHeaderMap"
$new.Range("B8").Value = "java.lang.NullPointerException"
$new.Range("B9").Value = "IN_MEMORY"
$new.Range("B10").Value = "INFO"

$new.Range("B11").Value = $false
$new.Range("B12").Value = $false
$new.Range("B13").Value = $false
$new.Range("B15").Value = $false
$new.Range("B17").Value = $false

$new.Range("B14").Value = 16
$new.Range("B16").Value = 0

# --- MEASUREMENTS block ---
$new.Range("B20").Value = 5118
$new.Range("B21").Value = 88079
$new.Range("B22").Value = 2917
$new.Range("B23").Value = 1964
$new.Range("B24").Value = 15468
$new.Range("B25").Value = 737
$new.Range("B26").Value = 469092
$new.Range("B27").Value = 69511
$new.Range("B28").Value = 902

# D20 formula needs updating to cover all 16 DDMIN rows
$new.Range("D20").Formula = "=COUNTA(A32:A47)"

# --- DDMIN data table rows 32-47 ---
$new.Range("A32").Value = "0-0"
$new.Range("B32").Value = 10
$new.Range("C32").Value = 10
$new.Range("D32").Formula = "=B32-C32"
$new.Range("E32").Value = 5118
$new.Range("F32").Value = 5118
$new.Range("G32").Formula = "=E32-F32"
$new.Range("H32").Value = 26
$new.Range("I32").Value = 1
$new.Range("J32").Formula = "=SUM($H$32:H32)"
$new.Range("K32").Formula = "=SUM($I$32:I32)"
$new.Range("L32").Value = 3808

$new.Range("A33").Value = "0-1"
$new.Range("B33").Value = 83
$new.Range("C33").Value = 83
$new.Range("D33").Formula = "=B33-C33"
$new.Range("E33").Value = 5118
$new.Range("F33").Value = 5118
$new.Range("G33").Formula = "=E33-F33"
$new.Range("H33").Value = 209
$new.Range("I33").Value = 0
$new.Range("J33").Formula = "=SUM($H$32:H33)"
$new.Range("K33").Formula = "=SUM($I$32:I33)"
$new.Range("L33").Value = 7450

$new.Range("A34").Value = "0-2"
$new.Range("B34").Value = 639
$new.Range("C34").Value = 550
$new.Range("D34").Formula = "=B34-C34"
$new.Range("E34").Value = 5118
$new.Range("F34").Value = 4248
$new.Range("G34").Formula = "=E34-F34"
$new.Range("H34").Value = 15861
$new.Range("I34").Value = 89
$new.Range("J34").Formula = "=SUM($H$32:H34)"
$new.Range("K34").Formula = "=SUM($I$32:I34)"
$new.Range("L34").Value = 217958

$new.Range("A35").Value = "0-3"
$new.Range("B35").Value = 795
$new.Range("C35").Value = 603
$new.Range("D35").Formula = "=B35-C35"
$new.Range("E35").Value = 4248
$new.Range("F35").Value = 4056
$new.Range("G35").Formula = "=E35-F35"
$new.Range("H35").Value = 41920
$new.Range("I35").Value = 598
$new.Range("J35").Formula = "=SUM($H$32:H35)"
$new.Range("K35").Formula = "=SUM($I$32:I35)"
$new.Range("L35").Value = 130337

$new.Range("A36").Value = "0-4"
$new.Range("B36").Value = 339
$new.Range("C36").Value = 169
$new.Range("D36").Formula = "=B36-C36"
$new.Range("E36").Value = 4056
$new.Range("F36").Value = 2248
$new.Range("G36").Formula = "=E36-F36"
$new.Range("H36").Value = 6397
$new.Range("I36").Value = 144
$new.Range("J36").Formula = "=SUM($H$32:H36)"
$new.Range("K36").Formula = "=SUM($I$32:I36)"
$new.Range("L36").Value = 65906

$new.Range("A37").Value = "0-5"
$new.Range("B37").Value = 154
$new.Range("C37").Value = 134
$new.Range("D37").Formula = "=B37-C37"
$new.Range("E37").Value = 2248
$new.Range("F37").Value = 2202
$new.Range("G37").Formula = "=E37-F37"
$new.Range("H37").Value = 1409
$new.Range("I37").Value = 20
$new.Range("J37").Formula = "=SUM($H$32:H37)"
$new.Range("K37").Formula = "=SUM($I$32:I37)"
$new.Range("L37").Value = 7650

$new.Range("A38").Value = "0-6"
$new.Range("B38").Value = 256
$new.Range("C38").Value = 242
$new.Range("D38").Formula = "=B38-C38"
$new.Range("E38").Value = 2202
$new.Range("F38").Value = 2038
$new.Range("G38").Formula = "=E38-F38"
$new.Range("H38").Value = 2733
$new.Range("I38").Value = 24
$new.Range("J38").Formula = "=SUM($H$32:H38)"
$new.Range("K38").Formula = "=SUM($I$32:I38)"
$new.Range("L38").Value = 27912

$new.Range("A39").Value = "0-7"
$new.Range("B39").Value = 85
$new.Range("C39").Value = 81
$new.Range("D39").Formula = "=B39-C39"
$new.Range("E39").Value = 2038
$new.Range("F39").Value = 2027
$new.Range("G39").Formula = "=E39-F39"
$new.Range("H39").Value = 505
$new.Range("I39").Value = 6
$new.Range("J39").Formula = "=SUM($H$32:H39)"
$new.Range("K39").Formula = "=SUM($I$32:I39)"
$new.Range("L39").Value = 2536

$new.Range("A40").Value = "0-8"
$new.Range("B40").Value = 44
$new.Range("C40").Value = 42
$new.Range("D40").Formula = "=B40-C40"
$new.Range("E40").Value = 2027
$new.Range("F40").Value = 2025
$new.Range("G40").Formula = "=E40-F40"
$new.Range("H40").Value = 191
$new.Range("I40").Value = 3
$new.Range("J40").Formula = "=SUM($H$32:H40)"
$new.Range("K40").Formula = "=SUM($I$32:I40)"
$new.Range("L40").Value = 1032

$new.Range("A41").Value = "0-9"
$new.Range("B41").Value = 30
$new.Range("C41").Value = 24
$new.Range("D41").Formula = "=B41-C41"
$new.Range("E41").Value = 2025
$new.Range("F41").Value = 1971
$new.Range("G41").Formula = "=E41-F41"
$new.Range("H41").Value = 175
$new.Range("I41").Value = 10
$new.Range("J41").Formula = "=SUM($H$32:H41)"
$new.Range("K41").Formula = "=SUM($I$32:I41)"
$new.Range("L41").Value = 1531

$new.Range("A42").Value = "0-10"
$new.Range("B42").Value = 13
$new.Range("C42").Value = 10
$new.Range("D42").Formula = "=B42-C42"
$new.Range("E42").Value = 1971
$new.Range("F42").Value = 1968
$new.Range("G42").Formula = "=E42-F42"
$new.Range("H42").Value = 37
$new.Range("I42").Value = 3
$new.Range("J42").Formula = "=SUM($H$32:H42)"
$new.Range("K42").Formula = "=SUM($I$32:I42)"
$new.Range("L42").Value = 588

$new.Range("A43").Value = "0-11"
$new.Range("B43").Value = 12
$new.Range("C43").Value = 10
$new.Range("D43").Formula = "=B43-C43"
$new.Range("E43").Value = 1968
$new.Range("F43").Value = 1966
$new.Range("G43").Formula = "=E43-F43"
$new.Range("H43").Value = 37
$new.Range("I43").Value = 3
$new.Range("J43").Formula = "=SUM($H$32:H43)"
$new.Range("K43").Formula = "=SUM($I$32:I43)"
$new.Range("L43").Value = 528

$new.Range("A44").Value = "0-12"
$new.Range("B44").Value = 2
$new.Range("C44").Value = 2
$new.Range("D44").Formula = "=B44-C44"
$new.Range("E44").Value = 1966
$new.Range("F44").Value = 1966
$new.Range("G44").Formula = "=E44-F44"
$new.Range("H44").Value = 2
$new.Range("I44").Value = 0
$new.Range("J44").Formula = "=SUM($H$32:H44)"
$new.Range("K44").Formula = "=SUM($I$32:I44)"
$new.Range("L44").Value = 120

$new.Range("A45").Value = "0-13"
$new.Range("B45").Value = 3
$new.Range("C45").Value = 2
$new.Range("D45").Formula = "=B45-C45"
$new.Range("E45").Value = 1966
$new.Range("F45").Value = 1964
$new.Range("G45").Formula = "=E45-F45"
$new.Range("H45").Value = 7
$new.Range("I45").Value = 1
$new.Range("J45").Formula = "=SUM($H$32:H45)"
$new.Range("K45").Formula = "=SUM($I$32:I45)"
$new.Range("L45").Value = 163

$new.Range("A46").Value = "0-14"
$new.Range("B46").Value = 2
$new.Range("C46").Value = 2
$new.Range("D46").Formula = "=B46-C46"
$new.Range("E46").Value = 1964
$new.Range("F46").Value = 1964
$new.Range("G46").Formula = "=E46-F46"
$new.Range("H46").Value = 2
$new.Range("I46").Value = 0
$new.Range("J46").Formula = "=SUM($H$32:H46)"
$new.Range("K46").Formula = "=SUM($I$32:I46)"
$new.Range("L46").Value = 22

$new.Range("A47").Value = "0-15"
$new.Range("B47").Value = 0
$new.Range("C47").Value = 0
$new.Range("D47").Formula = "=B47-C47"
$new.Range("E47").Value = 1964
$new.Range("F47").Value = 1964
$new.Range("G47").Formula = "=E47-F47"
$new.Range("H47").Value = 0
$new.Range("I47").Value = 0
$new.Range("J47").Formula = "=SUM($H$32:H47)"
$new.Range("K47").Formula = "=SUM($I$32:I47)"
$new.Range("L47").Value = 4

$tmpl.Activate()
